$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text so numeric-looking strings
# (e.g. "1.000", "229.10") are preserved exactly as typed, without
# Excel's automatic number coercion dropping trailing zeros.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range("D2").Value = "29.344.57"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "1.841.06"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("E5").Value = "  -0.64%  "

$ws.Range("D6").Value = "0.6282"
$ws.Range("E6").Value = "  -0.60%  "

$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "0.07438"
$ws.Range("E8").Value = "  -0.57%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.2892"
$ws.Range("E9").Value = "  -0.42%  "

$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "24.94"
$ws.Range("E10").Value = "  +2.34%  "

$ws.Range("D12").Value = "1.840.13"
$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").Value = "4.965"
$ws.Range("E13").Value = "  -0.67%  "

$ws.Range("D14").Value = "0.6755"
$ws.Range("E14").Value = "  -0.35%  "

$ws.Range("D15").Value = "0.00001027"
$ws.Range("E15").Value = "  +0.71%  "

$ws.Range("D16").Value = "81.56"
$ws.Range("E16").Value = "  -0.70%  "

$ws.Range("D17").Value = "6.236"
$ws.Range("E17").Value = "  +1.60%  "

$ws.Range("D18").Value = "29.392.42"
$ws.Range("E18").Value = "  +0.06%  "

$ws.Range("D19").Value = "229.10"
$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("E20").Value = "  -0.16%  "

$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").Value = "7.348"
$ws.Range("E22").Value = "  -1.24%  "

$ws.Range("D23").Value = "1.000"

$ws.Range("D24").Value = "158.15"
$ws.Range("E24").Value = "  -0.39%  "

$ws.Range("D25").Value = "8.475"
$ws.Range("E25").Value = "  +0.76%  "

$ws.Range("E26").Value = "  -2.41%  "

$ws.Range("D27").Value = "17.39"
$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("D28").Value = "0.07208"
$ws.Range("E28").Value = "  +14.14%  "

$ws.Range("E29").Value = "  +5.67%  "

$ws.Range("D30").Value = "1.480"
$ws.Range("E30").Value = "  +0.49%  "

$ws.Range("D31").Value = "4.047"
$ws.Range("E31").Value = "  -1.06%  "

$ws.Range("D32").Value = "4.042"
$ws.Range("E32").Value = "  -0.38%  "

$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").Value = "0.6962"

$ws.Range("D36").Value = "2.577"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").Value = "0.01842"
$ws.Range("E37").Value = "  +1.30%  "

$ws.Range("D38").Value = "2.801"
$ws.Range("E38").Value = "  -1.20%  "

$ws.Range("D39").Value = "1.235.61"
$ws.Range("E39").Value = "  -1.49%  "

$ws.Range("D40").Value = "6.809"
$ws.Range("E40").Value = "  +3.67%  "

$ws.Range("D41").Value = "0.9274"
$ws.Range("E41").Value = "  +2.22%  "

$ws.Range("E42").Value = "  +0.22%  "

$ws.Range("D43").Value = "2.004.37"
$ws.Range("E43").Value = "  -0.17%  "

$ws.Range("D44").Value = "100.41"
$ws.Range("E44").Value = "  -0.88%  "

$ws.Range("E45").Value = "  -1.49%  "

$ws.Range("E46").Value = "  +3.57%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.030"
$ws.Range("E47").Value = "  -0.25%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.712"
$ws.Range("E48").Value = "  +1.40%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.897"
$ws.Range("E49").Value = "  -1.60%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.1138"
$ws.Range("E50").Value = "  -3.09%  "

$ws.Range("D51").Value = "0.3907"
$ws.Range("E51").Value = "  -0.61%  "

# Restore the original (default/General) formatting on column D now
# that the text values are committed -- ClearFormats keeps the cell
# content as text (it does not re-run number inference) while
# dropping the temporary "@" number format we applied above, so the
# cells end up with no explicit style again, matching the source file.
$colD.ClearFormats()
